$d = $word.ActiveDocument

# Locate the target sentence: "Я, как бухгалтер, хочу получать информацию
# об оплаченных заказах в виде таблицы." and find the trailing period so we
# can split it into three separate runs:
#   1) "...в виде таблицы"
#   2) ", которую можно копировать куда-либо"
#   3) "."
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$found = $searchRange.Find.Execute(
    "бухгалтер, хочу получать информацию об оплаченных заказах в виде таблицы.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # The found range ends right after the trailing period (End points to the
    # paragraph mark that follows it), so the period itself is the very last
    # character of the range.
    $periodStart = $searchRange.End - 1
    $periodEnd = $searchRange.End

    $periodRange = $d.Range($periodStart, $periodEnd)

    $insertedXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>, которую можно копировать куда-либо</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $periodRange.InsertXML($insertedXml)
}
